$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24 - Initial rate (mol/s)
$ws.Range("A24").Value = "Initial rate (mol/s)"
$ws.Range("B24").Formula = "=(B22/1000000)*(B2/2)^2*PI()*B3*(1/1000)"
$ws.Range("C24").Formula = "=(C22/1000000)*(C2/2)^2*PI()*C3*(1/1000)"
$ws.Range("B24").NumberFormat = "0.00E+00"
$ws.Range("C24").NumberFormat = "0.00E+00"

# Row 25 - Initial rate (molecules/s)
$ws.Range("A25").Value = "Initial rate (molecules/s)"
$ws.Range("B25").Formula = "=B24*6.022E+23"
$ws.Range("C25").Formula = "=C24*6.022E+23"
$ws.Range("B25").NumberFormat = "0.00E+00"
$ws.Range("C25").NumberFormat = "0.00E+00"

# Row 26 - Quantum Yield
$ws.Range("A26").Value = "Quantum Yield"
$ws.Range("B26").Formula = "=(B25/B19)*100"
$ws.Range("C26").Formula = "=(C25/C19)*100"
$ws.Range("B26").NumberFormat = "0.00E+00"
$ws.Range("C26").NumberFormat = "0.00E+00"

$null = $ws.Range("C26").Select()
